$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8510499000549316
$ws.Range("B1").Value = 1.833649873733521
$ws.Range("D1").Value = 1.928653120994568
$ws.Range("E1").Value = 1.140200138092041
